$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 384, pushing the existing rows 384-403 down to 385-404.
$ws.Rows.Item(384).Insert()

# Populate the newly inserted row 384 with the new data point
$ws.Range("A384").Value = 11
$ws.Range("B384").Value = "Vega Monumental Concepción"
$ws.Range("C384").Value = "Bíobío"
$ws.Range("D384").Value = 45041
$ws.Range("E384").Value = 8
$ws.Range("F384").Value = 100112009
$ws.Range("G384").Value = "Acelga"
$ws.Range("H384").Value = "Sin especificar"
$ws.Range("I384").Value = "Primera"
$ws.Range("J384").Value = 400
$ws.Range("K384").Value = 600
$ws.Range("L384").Value = 650
$ws.Range("M384").Value = 625
$ws.Range("N384").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O384").Value = "Región Metropolitana"
$ws.Range("P384").Value = 625
$ws.Range("Q384").Value = 1
$ws.Range("R384").Value = "Hortaliza"
